$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2150
$ws.Range("I40").Value = 2150
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2150
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1975
$ws.Range("N40").ClearContents()
$ws.Range("H51").Value = 12000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 12000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 12000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -12968
$ws.Range("H55").Value = 373.2
$ws.Range("I55").Value = 272
$ws.Range("J55").Value = 525
$ws.Range("K55").Value = 272
$ws.Range("L55").Value = 525
$ws.Range("M55").Value = -58
$ws.Range("N55").Value = -953
$ws.Range("H113").Value = 7895.857
$ws.Range("I113").Value = 7544.75
$ws.Range("K113").Value = 7544.75
$ws.Range("M113").Value = -4290.75
$ws.Range("H132").Value = 2201.9565
$ws.Range("I132").Value = 1924.7778
$ws.Range("K132").Value = 5774.3334
$ws.Range("M132").Value = -3244.3334
$ws.Range("H135").Value = 1799.5
$ws.Range("I135").Value = 429.42856
$ws.Range("K135").Value = 3864.85704
$ws.Range("M135").Value = -1329.85704
$ws.Range("H137").Value = 1125.1428
$ws.Range("I137").Value = 1125.1428
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3375.4284
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -825.4284000000002
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 4133
$ws.Range("I138").Value = 4439.2
$ws.Range("J138").Value = 3914.2856
$ws.Range("K138").Value = 13317.6
$ws.Range("L138").Value = 11742.8568
$ws.Range("M138").Value = -8177.599999999999
$ws.Range("N138").Value = -22022.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4669.4287
$ws.Range("I32").Value = 4414.2
$ws.Range("K32").Value = 4414.2
$ws.Range("M32").Value = -4127.2
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H80").Value = 40550
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -61996
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H83").Value = 40550
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -189984
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H110").Value = 2911
$ws.Range("I110").Value = 2911
$ws.Range("K110").Value = 2911
$ws.Range("M110").Value = -866
$ws.Range("H132").Value = 916.75
$ws.Range("I132").Value = 916.75
$ws.Range("K132").Value = 2750.25
$ws.Range("M132").Value = -220.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1908.9
$ws.Range("I31").Value = 1028
$ws.Range("J31").Value = 3964.3333
$ws.Range("K31").Value = 1028
$ws.Range("L31").Value = 3964.3333
$ws.Range("M31").Value = -733
$ws.Range("N31").Value = -4554.3333
$ws.Range("H34").Value = 1908.9
$ws.Range("I34").Value = 1028
$ws.Range("J34").Value = 3964.3333
$ws.Range("K34").Value = 1028
$ws.Range("L34").Value = 3964.3333
$ws.Range("M34").Value = -826
$ws.Range("N34").Value = -4368.3333
$ws.Range("H122").Value = 1880.7693
$ws.Range("I122").Value = 1926.7
$ws.Range("K122").Value = 5780.1
$ws.Range("M122").Value = -3330.1
$ws.Range("H134").Value = 798.2857
$ws.Range("I134").Value = 798.2857
$ws.Range("K134").Value = 2394.8571
$ws.Range("M134").Value = 140.1428999999998
$ws.Range("H141").Value = 122897.414
$ws.Range("I141").Value = 59000
$ws.Range("J141").Value = 128706.27
$ws.Range("K141").Value = 59000
$ws.Range("L141").Value = 128706.27
$ws.Range("M141").Value = -53820
$ws.Range("N141").Value = -139066.27

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 883.3333
$ws.Range("I5").Value = 825.5
$ws.Range("K5").Value = 2476.5
$ws.Range("M5").Value = -2364.5
$ws.Range("H34").Value = 634.3
$ws.Range("I34").Value = 655.375
$ws.Range("J34").Value = 550
$ws.Range("K34").Value = 1966.125
$ws.Range("L34").Value = 1650
$ws.Range("M34").Value = -1882.125
$ws.Range("N34").Value = -1818
$ws.Range("H39").Value = 2066.6667
$ws.Range("J39").Value = 2066.6667
$ws.Range("L39").Value = 6200.000100000001
$ws.Range("N39").Value = -6788.000100000001
$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H69").Value = 3157.6
$ws.Range("I69").Value = 598
$ws.Range("J69").Value = 6997
$ws.Range("K69").Value = 1794
$ws.Range("L69").Value = 20991
$ws.Range("M69").Value = -983
$ws.Range("N69").Value = -22613
$ws.Range("H72").Value = 3157.6
$ws.Range("I72").Value = 598
$ws.Range("J72").Value = 6997
$ws.Range("K72").Value = 5382
$ws.Range("L72").Value = 62973
$ws.Range("M72").Value = -1326
$ws.Range("N72").Value = -71085
$ws.Range("H80").Value = 5414.1665
$ws.Range("I80").Value = 5998.75
$ws.Range("J80").Value = 5121.875
$ws.Range("K80").Value = 17996.25
$ws.Range("L80").Value = 15365.625
$ws.Range("M80").Value = -17060.25
$ws.Range("N80").Value = -17237.625
$ws.Range("H83").Value = 5414.1665
$ws.Range("I83").Value = 5998.75
$ws.Range("J83").Value = 5121.875
$ws.Range("K83").Value = 53988.75
$ws.Range("L83").Value = 46096.875
$ws.Range("M83").Value = -49308.75
$ws.Range("N83").Value = -55456.875
$ws.Range("H117").Value = 3118.875
$ws.Range("I117").Value = 2500
$ws.Range("J117").Value = 3207.2856
$ws.Range("K117").Value = 7500
$ws.Range("L117").Value = 9621.856800000001
$ws.Range("M117").Value = -4058
$ws.Range("N117").Value = -16505.8568
$ws.Range("H129").Value = 2489.4285
$ws.Range("I129").Value = 1276
$ws.Range("J129").Value = 3399.5
$ws.Range("K129").Value = 3828
$ws.Range("L129").Value = 10198.5
$ws.Range("M129").Value = 1172
$ws.Range("N129").Value = -20198.5
$ws.Range("H131").Value = 892
$ws.Range("I131").Value = 913.4
$ws.Range("K131").Value = 2740.2
$ws.Range("M131").Value = 2299.8
$ws.Range("H135").Value = 883.3333
$ws.Range("I135").Value = 825.5
$ws.Range("K135").Value = 7429.5
$ws.Range("M135").Value = -4894.5
$ws.Range("H137").Value = 3422.5
$ws.Range("J137").Value = 3585.682
$ws.Range("L137").Value = 10757.046
$ws.Range("N137").Value = -20957.046

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2243
$ws.Range("I122").Value = 1907.6666
$ws.Range("J122").Value = 3249
$ws.Range("K122").Value = 5722.9998
$ws.Range("L122").Value = 9747
$ws.Range("M122").Value = -3272.9998
$ws.Range("N122").Value = -14647

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5006.8
$ws.Range("I132").Value = 5006.8
$ws.Range("K132").Value = 15020.4
$ws.Range("M132").Value = -12490.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H113").Value = 1945
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 1490
$ws.Range("K113").Value = 7200
$ws.Range("L113").Value = 4470
$ws.Range("M113").Value = -5030
$ws.Range("N113").Value = -8810
$ws.Range("H132").Value = 1917.909
$ws.Range("I132").Value = 1630.875
$ws.Range("J132").Value = 2683.3333
$ws.Range("K132").Value = 4892.625
$ws.Range("L132").Value = 8049.999899999999
$ws.Range("M132").Value = -2362.625
$ws.Range("N132").Value = -13109.9999
$ws.Range("H136").Value = 1679.9231
$ws.Range("I136").Value = 1508.9
$ws.Range("K136").Value = 4526.700000000001
$ws.Range("M136").Value = -1976.700000000001

